# Correction in SA algorithm and 746 logs
# Updates the "Fitness" column (C) values for run_18 log rows so that the
# recorded best-fitness-so-far sequence reflects the corrected algorithm.
# Data rows 2-190 (Generation 0-188) get corrected values; rows 191-252
# already hold the correct values and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is (startRow, endRow, newFitnessValue)
$updates = @(
    @(2, 4, 8477),
    @(5, 6, 8288),
    @(7, 8, 7817),
    @(9, 16, 7667),
    @(17, 71, 7312),
    @(92, 190, 7293)
)

foreach ($u in $updates) {
    $startRow = $u[0]
    $endRow = $u[1]
    $value = $u[2]
    $rangeAddr = "C$startRow`:C$endRow"
    $ws.Range($rangeAddr).Value = $value
}
